$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.945.05'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '2.261.93'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.91'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.23'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.526'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.95%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.05'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +6.92%  '
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.66'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').Value = '2.610.05'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.40'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').Value = '2.270.91'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.792'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').Value = '41.833.53'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.37'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('D20').Value = '0.0₃0901'
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.14'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.93'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('E24').Value = '  -1.64%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.93'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.52'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.46'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.11'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '160.12'
$ws.Range('D31').ClearFormats()
$ws.Range('E32').Value = '  -2.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.17'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0735'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.98'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('E40').Value = '  -2.53%  '
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('E42').Value = '  +3.11%  '
$ws.Range('D43').Value = '1.966.61'
$ws.Range('E43').Value = '  -2.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0283'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.81'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.56%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.85'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.14'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '72.70'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '91.10'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.02%  '
